# GameData.xlsx update: simple monster-tracking AI via Behavior Designer,
# Navmesh AI package added. This adjusts the SkillData sheet: a few skills
# are reclassified from Active to Passive (with Debuff buff type and
# Hp/Mp cost types), skill ranges/targeting counts are rebalanced, and the
# active-cell selection on that sheet moves to M13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SkillData")
$ws.Activate() | Out-Null

# Row 2
$ws.Range("J2").Value = 100
$ws.Range("O2").Value = 2

# Row 3
$ws.Range("K3").Value = 11

# Row 4
$ws.Range("D4").Value = "Passive"
$ws.Range("H4").Value = "Debuff"
$ws.Range("I4").Value = "Hp"
$ws.Range("K4").Value = 13
$ws.Range("O4").Value = 4

# Row 5
$ws.Range("D5").Value = "Passive"
$ws.Range("H5").Value = "Debuff"
$ws.Range("I5").Value = "Mp"
$ws.Range("K5").Value = 3
$ws.Range("N5").Value = "AutoTargeting"
$ws.Range("O5").Value = 5

# Row 6
$ws.Range("D6").Value = "Passive"
$ws.Range("H6").Value = "Debuff"
$ws.Range("K6").Value = 6
$ws.Range("O6").Value = 6

# Row 7
$ws.Range("D7").Value = "Passive"
$ws.Range("O7").Value = 7

# Row 8
$ws.Range("D8").Value = "Passive"
$ws.Range("O8").Value = 8

# Selection moves to M13 on the SkillData sheet
$ws.Range("M13").Select() | Out-Null
